# Update Data Sources from LFX 2026-02-24 (#98)
#
# The upstream sync replaced the table style applied to every data-source
# table in the deck (old style {928D4016-75EB-4100-83A0-C58B1C6B50E9} ->
# new style {2148295F-C88A-4D38-9292-040BEF2578C8}).
#
# Walk every slide/shape, find the tables, and re-apply the new table
# style via Table.ApplyStyle (Table.Style is read-only in this object
# model and must be changed through the ApplyStyle method).

$oldStyleId = "{928D4016-75EB-4100-83A0-C58B1C6B50E9}"
$newStyleId = "{2148295F-C88A-4D38-9292-040BEF2578C8}"

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
